$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Duplicate row 5 (formatting only) down into the two new rows, 6 and 7 ---
$ws.Range("A5:N5").Copy()
$ws.Range("A6:N6").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A5:N5").Copy()
$ws.Range("A7:N7").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Match row 5's explicit custom row height behaviour (14.25pt, customHeight)
$ws.Rows.Item(6).RowHeight = 14.25
$ws.Rows.Item(7).RowHeight = 14.25

# --- Row 6 values: Tashulya / Hamidova ---
$ws.Cells.Item(6,1).Value = "Tashulya"
$ws.Cells.Item(6,2).Value = "Hamidova"
$ws.Cells.Item(6,3).Value = 45940
$ws.Cells.Item(6,4).Value = 0.41666666666666702
$ws.Cells.Item(6,5).Value = 45945
$ws.Cells.Item(6,6).Value = 0.375
$ws.Cells.Item(6,7).Value = 45975
$ws.Cells.Item(6,8).Value = 0.58333333333333337
$ws.Cells.Item(6,9).Value = 45977
$ws.Cells.Item(6,10).Value = 0.58333333333333304
$ws.Cells.Item(6,11).Value = 46005
$ws.Cells.Item(6,12).Value = 0.625
$ws.Cells.Item(6,13).Value = 46007
$ws.Cells.Item(6,14).Value = 0.70833333333333337

# --- Row 7 values: same names / schedule as row 6 ---
$ws.Cells.Item(7,1).Value = "Tashulya"
$ws.Cells.Item(7,2).Value = "Hamidova"
$ws.Cells.Item(7,3).Value = 45940
$ws.Cells.Item(7,4).Value = 0.41666666666666702
$ws.Cells.Item(7,5).Value = 45945
$ws.Cells.Item(7,6).Value = 0.375
$ws.Cells.Item(7,7).Value = 45975
$ws.Cells.Item(7,8).Value = 0.58333333333333337
$ws.Cells.Item(7,9).Value = 45977
$ws.Cells.Item(7,10).Value = 0.58333333333333304
$ws.Cells.Item(7,11).Value = 46005
$ws.Cells.Item(7,12).Value = 0.625
$ws.Cells.Item(7,13).Value = 46007
$ws.Cells.Item(7,14).Value = 0.70833333333333337

# Re-apply the date number format on the new rows' date cells so the workbook
# carries its own (duplicated) date style, matching the freshly-added style
# record used for the October/November/December date columns in the new rows.
$ws.Cells.Item(6,3).NumberFormat = "mm-dd-yy"
$ws.Cells.Item(6,5).NumberFormat = "mm-dd-yy"
$ws.Cells.Item(6,7).NumberFormat = "mm-dd-yy"
$ws.Cells.Item(6,9).NumberFormat = "mm-dd-yy"
$ws.Cells.Item(6,11).NumberFormat = "mm-dd-yy"
$ws.Cells.Item(6,13).NumberFormat = "mm-dd-yy"
$ws.Cells.Item(7,3).NumberFormat = "mm-dd-yy"
$ws.Cells.Item(7,5).NumberFormat = "mm-dd-yy"
$ws.Cells.Item(7,7).NumberFormat = "mm-dd-yy"
$ws.Cells.Item(7,9).NumberFormat = "mm-dd-yy"
$ws.Cells.Item(7,11).NumberFormat = "mm-dd-yy"
$ws.Cells.Item(7,13).NumberFormat = "mm-dd-yy"
